# The sheet is a daily log of price records (one row per observation),
# ordered from most recent date to oldest. A new observation was added
# at the top of the data block (row 104), pushing all subsequent rows
# down by one. This grows the used range from A1:R143 to A1:R144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104; this shifts rows 104:143 down to
# 105:144 (carrying their values and formatting, e.g. the date style
# on column D) and automatically extends the sheet dimension.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new observation.
$ws.Cells.Item(104, 1).Value  = 3
$ws.Cells.Item(104, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(104, 3).Value  = "Coquimbo"
$ws.Cells.Item(104, 4).Value  = 44726
$ws.Cells.Item(104, 5).Value  = 5
$ws.Cells.Item(104, 6).Value  = 100112026
$ws.Cells.Item(104, 7).Value  = "Haba"
$ws.Cells.Item(104, 8).Value  = "Sin especificar"
$ws.Cells.Item(104, 9).Value  = "Primera"
$ws.Cells.Item(104, 10).Value = 82
$ws.Cells.Item(104, 11).Value = 22000
$ws.Cells.Item(104, 12).Value = 23000
$ws.Cells.Item(104, 13).Value = 22512
$ws.Cells.Item(104, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(104, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(104, 16).Value = 900
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = "Hortaliza"
